$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.328.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.668.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.70%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'220.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.95%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06367"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07837"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.673.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.897.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.95%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'65.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.60%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.337.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.724"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.41%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'198.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +1.67%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.058"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.51%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.009"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.15%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1219"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.253"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.88%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.507"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.44%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05905"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.555"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.19%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.65%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.98%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.434"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5822"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.72%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01616"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.80%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.958"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'1.076.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.10%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8593"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.87%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.02%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'102.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.807.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.59%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'58.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.12%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.014"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.64%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4405"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.82%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.55%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.075"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.20%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.11%  "
$ws.Range("E51").Style = "Normal"
